$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 7.5
$ws.Range("N2").Value = 1.62
$ws.Range("O2").Value = 2.3
$ws.Range("X2").Value = 12
$ws.Range("AA2").Value = 10
$ws.Range("AD2").Value = 351
$ws.Range("AE2").Value = 21
$ws.Range("AG2").Value = 23
$ws.Range("AH2").Value = 101

# Row 3
$ws.Range("G3").Value = 2.3
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 3
$ws.Range("J3").Value = 1.05
$ws.Range("K3").Value = 11
$ws.Range("N3").Value = 1.88
$ws.Range("O3").Value = 2.02
$ws.Range("P3").Value = 1.36
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = 1.67
$ws.Range("S3").Value = 2.1
$ws.Range("W3").Value = 26
$ws.Range("X3").Value = 21
$ws.Range("Y3").Value = 29
$ws.Range("Z3").Value = 12
$ws.Range("AA3").Value = 7
$ws.Range("AB3").Value = 15
$ws.Range("AE3").Value = 12
$ws.Range("AF3").Value = 17
$ws.Range("AJ3").Value = 34

# Row 4
$ws.Range("G4").Value = 1.2
$ws.Range("H4").Value = 7.5
$ws.Range("I4").Value = 12
$ws.Range("R4").Value = 1.8
$ws.Range("S4").Value = 1.95
$ws.Range("T4").Value = 12
$ws.Range("U4").Value = 8.5
$ws.Range("V4").Value = 10
$ws.Range("W4").Value = 9
$ws.Range("Y4").Value = 26
$ws.Range("AA4").Value = 15
$ws.Range("AB4").Value = 26
$ws.Range("AC4").Value = 67
$ws.Range("AD4").Value = 301
$ws.Range("AE4").Value = 41

# Row 5
$ws.Range("G5").Value = 1.65
$ws.Range("I5").Value = 4.75
$ws.Range("J5").Value = 1.03
$ws.Range("K5").Value = 17
$ws.Range("L5").Value = 1.17
$ws.Range("M5").Value = 5
$ws.Range("N5").Value = 1.53
$ws.Range("O5").Value = 2.5
$ws.Range("P5").Value = 1.29
$ws.Range("Q5").Value = 3.5
$ws.Range("R5").Value = 1.57
$ws.Range("S5").Value = 2.25
$ws.Range("T5").Value = 11
$ws.Range("U5").Value = 11
$ws.Range("V5").Value = 9
$ws.Range("W5").Value = 15
$ws.Range("Z5").Value = 17
$ws.Range("AA5").Value = 9
$ws.Range("AB5").Value = 15
$ws.Range("AD5").Value = 151
$ws.Range("AE5").Value = 19
$ws.Range("AG5").Value = 17
$ws.Range("AI5").Value = 41
$ws.Range("AJ5").Value = 34

# Row 6
$ws.Range("G6").Value = 5.5
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 1.6
$ws.Range("J6").Value = 1.04
$ws.Range("K6").Value = 12
$ws.Range("N6").Value = 1.8
$ws.Range("O6").Value = 2
$ws.Range("R6").Value = 1.8
$ws.Range("S6").Value = 1.95
$ws.Range("V6").Value = 19
$ws.Range("X6").Value = 41
$ws.Range("Y6").Value = 41
$ws.Range("AD6").Value = 301
$ws.Range("AF6").Value = 8.5
$ws.Range("AH6").Value = 13

# Row 7
$ws.Range("G7").Value = 1.48
$ws.Range("I7").Value = 6.5
$ws.Range("L7").Value = 1.25
$ws.Range("M7").Value = 4
$ws.Range("N7").Value = 1.86
$ws.Range("O7").Value = 2.04
$ws.Range("R7").Value = 1.95
$ws.Range("S7").Value = 1.8
$ws.Range("T7").Value = 7.5
$ws.Range("X7").Value = 13
$ws.Range("Y7").Value = 29
$ws.Range("Z7").Value = 12
$ws.Range("AB7").Value = 21
$ws.Range("AC7").Value = 67
$ws.Range("AD7").Value = 451
$ws.Range("AF7").Value = 41

# Row 8
$ws.Range("G8").Value = 1.4
$ws.Range("I8").Value = 7.5
$ws.Range("R8").Value = 2.1
$ws.Range("S8").Value = 1.67
$ws.Range("U8").Value = 6
$ws.Range("W8").Value = 9
$ws.Range("X8").Value = 12
$ws.Range("AA8").Value = 9
$ws.Range("AD8").Value = 501
$ws.Range("AE8").Value = 17
$ws.Range("AF8").Value = 41

# Row 15
$ws.Range("G15").Value = 2.2
$ws.Range("H15").Value = 2.88
$ws.Range("J15").Value = 1.14
$ws.Range("K15").Value = 5.5
$ws.Range("L15").Value = 1.67
$ws.Range("M15").Value = 2.2
$ws.Range("N15").Value = 3.2
$ws.Range("O15").Value = 1.36
$ws.Range("P15").Value = 1.67
$ws.Range("Q15").Value = 2.1
$ws.Range("R15").Value = 2.38
$ws.Range("S15").Value = 1.53
$ws.Range("U15").Value = 8.5
$ws.Range("X15").Value = 23
$ws.Range("Z15").Value = 5.5
$ws.Range("AB15").Value = 21
$ws.Range("AE15").Value = 7.5

# Row 26
$ws.Range("H26").Value = 8.5
$ws.Range("R26").Value = 2.1
$ws.Range("S26").Value = 1.67
$ws.Range("T26").Value = 9.5
$ws.Range("U26").Value = 6.5
$ws.Range("AA26").Value = 17
$ws.Range("AC26").Value = 81
$ws.Range("AD26").Value = 501
$ws.Range("AE26").Value = 29

# Row 27
$ws.Range("G27").Value = 3.4
$ws.Range("H27").Value = 3.25
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = 1.08
$ws.Range("K27").Value = 7
$ws.Range("N27").Value = 2.3
$ws.Range("O27").Value = 1.6
$ws.Range("T27").Value = 9
$ws.Range("U27").Value = 17
$ws.Range("AB27").Value = 17
$ws.Range("AF27").Value = 9

# Row 28
$ws.Range("G28").Value = 2.3
$ws.Range("I28").Value = 2.6
$ws.Range("J28").Value = 1.01
$ws.Range("K28").Value = 13
$ws.Range("N28").Value = 1.73
$ws.Range("O28").Value = 2.08
$ws.Range("R28").Value = 1.62
$ws.Range("S28").Value = 2.2
$ws.Range("U28").Value = 13
$ws.Range("X28").Value = 19
$ws.Range("Y28").Value = 23
$ws.Range("AB28").Value = 13
$ws.Range("AH28").Value = 26

# Row 29
$ws.Range("N29").Value = 1.7
$ws.Range("O29").Value = 2.1
$ws.Range("T29").Value = 8
$ws.Range("AA29").Value = 8.5
$ws.Range("AD29").Value = 201

# Row 30
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 3.6
$ws.Range("I30").Value = 3.5
$ws.Range("J30").Value = 1.04
$ws.Range("K30").Value = 12
$ws.Range("L30").Value = 1.25
$ws.Range("M30").Value = 3.75
$ws.Range("N30").Value = 1.83
$ws.Range("O30").Value = 2.03
$ws.Range("U30").Value = 10
$ws.Range("W30").Value = 17
$ws.Range("X30").Value = 15
$ws.Range("AE30").Value = 12
$ws.Range("AF30").Value = 19
$ws.Range("AH30").Value = 41
$ws.Range("AJ30").Value = 34
